# Adds two new Q&A rows (43 and 44) to the "SQL" worksheet, mirroring the
# existing pattern used for rows further up the sheet: column A holds the
# (text-formatted) row number, column B the multi-line SQL query text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SQL")

# --- Row 43: "42" / PTNT_ACCT_NBR query --------------------------------
$query42 = @'
Select c.PTNT_ACCT_NBR, cp.SETL_DT
from PP001.CONSOLIDATED_PAYMENT cp, PP001.PROVIDER p, PP001.CLAIM c, PP001.UNCONSOLIDATED_PAYMENT ucp, PP001.CLAIM_UNCONSOLIDATED_PAYMENT cup
where cp.prov_key_id = p.prov_key_id
and cp.CONSL_PAY_NBR = ucp.CONSL_PAY_NBR
and ucp.UCONSL_PAY_KEY_ID = cup.UCONSL_PAY_KEY_ID
and cup.CLM_KEY_ID = c.CLM_KEY_ID
and p.PROV_TAX_ID_NBR = '{$tin}'
and cp.SETL_DT <= current date 
order by cp.SETL_DT DESC
fetch first row only

'@

# --- Row 44: "43" / SBSCR_ID query --------------------------------------
$query43 = @'
Select sr.SBSCR_ID, cp.SETL_DT
from PP001.CONSOLIDATED_PAYMENT cp, PP001.PROVIDER p, PP001.CLAIM c, PP001.UNCONSOLIDATED_PAYMENT ucp, PP001.CLAIM_UNCONSOLIDATED_PAYMENT cup, PP001.SUBSCRIBER sr
where cp.prov_key_id = p.prov_key_id
and cp.CONSL_PAY_NBR = ucp.CONSL_PAY_NBR
and ucp.UCONSL_PAY_KEY_ID = cup.UCONSL_PAY_KEY_ID
and cup.CLM_KEY_ID = c.CLM_KEY_ID
and c.SBSCR_KEY_ID = sr.SBSCR_KEY_ID
and p.PROV_TAX_ID_NBR = '840611484'
and cp.SETL_DT <= current date 
order by cp.SETL_DT DESC
fetch first row only
'@

# Seed the two new rows by copying the formatting (styles) of the last
# existing row (42: quote-prefixed "number-as-text" in A, wrap-text in B)
# so the new cells inherit identical styles (s="5" / s="7").
$ws.Range("A42:B42").Copy() | Out-Null
$ws.Range("A43:B44").PasteSpecial(-4122) | Out-Null

# Column A/B values are written in the same order the shared-strings table
# picks them up (42, query42, query43, 43) so new <si> entries land at the
# expected indices.  Column A values are stored as text (quote-prefixed),
# matching how the existing sheet stores its "row number" labels.
$ws.Range("A43").Value = "'42"
$ws.Range("B43").Value = $query42
$ws.Range("B44").Value = $query43
$ws.Range("A44").Value = "'43"

# Match the row heights Excel computed for these wrapped cells.
$ws.Rows.Item(43).RowHeight = 180
$ws.Rows.Item(44).RowHeight = 195

# Update the selection to reflect where Excel leaves the cursor after the
# edit (bottom of the newly entered data).
$ws.Range("C44").Select() | Out-Null
